# "Registro de un usuario completo"
# Adds a new row (row 8) to the page-titles sheet for the "Register" page:
#   A8 = "Register: Mercury Tours"  (wrapped text, like the other page-name cells)
#   B8 = "Registro"                (short title, like the other short-title cells)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the bottom so it inherits the formatting of row 7 above it
# (this gives A8/B8 the same base styles used by the existing data rows).
$ws.Rows("8").Insert()

$ws.Range("A8").Value = "Register: Mercury Tours"
$ws.Range("B8").Value = "Registro"

# The long page-title column wraps its text.
$ws.Range("A8").WrapText = $true

# Row grows a bit taller to fit the wrapped text.
$ws.Rows("8").RowHeight = 15.65

# Leave the selection on the last cell that was entered.
$ws.Range("B8").Select() | Out-Null
